$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.431.64'
$ws.Range('E2').Value = '  -1.46%  '
$ws.Range('D3').Value = '1.839.45'
$ws.Range('E3').Value = '  -2.09%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '261.02'
$ws.Range('E5').Value = '  -5.97%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5185'
$ws.Range('E7').Value = '  -1.92%  '
$ws.Range('E8').Value = '  -4.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06775'
$ws.Range('E9').Value = '  -2.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.62'
$ws.Range('E10').Value = '  -6.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7678'
$ws.Range('E11').Value = '  -4.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07708'
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').Value = '1.828.73'
$ws.Range('E13').Value = '  -2.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.21'
$ws.Range('E14').Value = '  -2.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.031'
$ws.Range('E15').Value = '  -2.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9999'
$ws.Range('E17').Value = '  -4.41%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007973'
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9999'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '26.445.91'
$ws.Range('E20').Value = '  -1.57%  '
$ws.Range('D21').Value = '2.068.94'
$ws.Range('E21').Value = '  -2.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.575'
$ws.Range('E22').Value = '  -3.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.468'
$ws.Range('E23').Value = '  -5.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.958'
$ws.Range('E24').Value = '  -3.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.32'
$ws.Range('E25').Value = '  -1.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.213'
$ws.Range('E26').Value = '  -7.10%  '
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.98'
$ws.Range('E28').Value = '  -2.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '111.27'
$ws.Range('E29').Value = '  -1.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.186'
$ws.Range('E30').Value = '  -3.53%  '
$ws.Range('E31').Value = '  -3.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08725'
$ws.Range('E32').Value = '  -1.79%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04806'
$ws.Range('E33').Value = '  -1.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.128'
$ws.Range('E34').Value = '  -3.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.837'
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7059'
$ws.Range('E36').Value = '  -2.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.079'
$ws.Range('E37').Value = '  -6.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.217'
$ws.Range('E38').Value = '  -5.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01759'
$ws.Range('E39').Value = '  -4.34%  '
$ws.Range('E40').Value = '  -5.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '111.34'
$ws.Range('E41').Value = '  -4.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8905'
$ws.Range('E42').Value = '  -6.75%  '
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9994'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.689'
$ws.Range('E45').Value = '  -4.80%  '
$ws.Range('E46').Value = '  -6.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05865'
$ws.Range('E47').Value = '  -1.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.964'
$ws.Range('E48').Value = '  -3.74%  '
$ws.Range('E49').Value = '  -3.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1219'
$ws.Range('E50').Value = '  -8.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8871'
$ws.Range('E51').Value = '  +0.66%  '
